$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Brazil, Government of'
$ws.Cells.Item(2, 2).Value = 113766
$ws.Cells.Item(3, 1).Value = 'ExxonMobil'
$ws.Cells.Item(3, 2).Value = 300000
$ws.Cells.Item(4, 1).Value = 'South32'
$ws.Cells.Item(4, 2).Value = 250000
$ws.Cells.Item(5, 1).Value = 'United Kingdom, Government of'
$ws.Cells.Item(5, 2).Value = 57002344
$ws.Cells.Item(6, 1).Value = 'United States of America, Government of'
$ws.Cells.Item(6, 2).Value = 138658722
$ws.Cells.Item(7, 1).Value = 'Portugal, Government of'
$ws.Cells.Item(7, 2).Value = 7897114
$ws.Cells.Item(8, 1).Value = 'Qatar, Government of'
$ws.Cells.Item(8, 2).Value = 3000000
$ws.Cells.Item(9, 1).Value = 'Denmark, Government of'
$ws.Cells.Item(9, 2).Value = 1938222
$ws.Cells.Item(10, 1).Value = 'Germany, Government of'
$ws.Cells.Item(10, 2).Value = 18131405
$ws.Cells.Item(11, 1).Value = 'Spain, Government of'
$ws.Cells.Item(11, 2).Value = 1281514
$ws.Cells.Item(12, 1).Value = 'Australia, Government of'
$ws.Cells.Item(12, 2).Value = 884643
$ws.Cells.Item(13, 1).Value = 'European Commission''s Humanitarian Aid and Civil Protection Department'
$ws.Cells.Item(13, 2).Value = 40660722
$ws.Cells.Item(14, 1).Value = 'France, Government of'
$ws.Cells.Item(14, 2).Value = 1116435
$ws.Cells.Item(15, 1).Value = 'Ireland, Government of'
$ws.Cells.Item(15, 2).Value = 3448558
$ws.Cells.Item(16, 1).Value = 'Swiss Solidarity'
$ws.Cells.Item(16, 2).Value = 4760662
$ws.Cells.Item(17, 1).Value = 'Disasters Emergency Committee Cyclone Idai'
$ws.Cells.Item(17, 2).Value = 10400740
$ws.Cells.Item(18, 1).Value = 'Korea, Republic of, Government of'
$ws.Cells.Item(18, 2).Value = 680000
$ws.Cells.Item(19, 1).Value = 'Samaritan''s Purse'
$ws.Cells.Item(19, 2).Value = 138485
$ws.Cells.Item(20, 1).Value = 'United Nations Children''s Fund'
$ws.Cells.Item(20, 2).Value = 23952171
$ws.Cells.Item(21, 1).Value = 'Austria, Government of'
$ws.Cells.Item(21, 2).Value = 4383176
$ws.Cells.Item(22, 1).Value = 'Finland, Government of'
$ws.Cells.Item(22, 2).Value = 1114893
$ws.Cells.Item(23, 1).Value = 'TechnipFMC'
$ws.Cells.Item(23, 2).Value = 100000
$ws.Cells.Item(24, 1).Value = 'Kuwait, Government of'
$ws.Cells.Item(24, 2).Value = 416490
$ws.Cells.Item(25, 1).Value = 'Norway, Government of'
$ws.Cells.Item(25, 2).Value = 9577423
$ws.Cells.Item(26, 1).Value = 'Japan, Government of'
$ws.Cells.Item(26, 2).Value = 12188030
$ws.Cells.Item(27, 1).Value = 'Switzerland, Government of'
$ws.Cells.Item(27, 2).Value = 3886420
$ws.Cells.Item(28, 1).Value = 'Canada, Government of'
$ws.Cells.Item(28, 2).Value = 8733690
$ws.Cells.Item(29, 1).Value = 'Sweden, Government of'
$ws.Cells.Item(29, 2).Value = 17176387
$ws.Cells.Item(30, 1).Value = 'Central Emergency Response Fund'
$ws.Cells.Item(30, 2).Value = 35944553
$ws.Cells.Item(31, 1).Value = 'UNICEF National Committee/Australia'
$ws.Cells.Item(31, 2).Value = 25948
$ws.Cells.Item(32, 1).Value = 'UNICEF National Committee/Denmark'
$ws.Cells.Item(32, 2).Value = 82889
$ws.Cells.Item(33, 1).Value = 'UNICEF National Committee/Canada'
$ws.Cells.Item(33, 2).Value = 43694
$ws.Cells.Item(34, 1).Value = 'UNICEF National Committee/Netherlands'
$ws.Cells.Item(34, 2).Value = 81901
$ws.Cells.Item(35, 1).Value = 'UNICEF National Committee/Spain'
$ws.Cells.Item(35, 2).Value = 594377
$ws.Cells.Item(36, 1).Value = 'UNICEF National Committee/Norway'
$ws.Cells.Item(36, 2).Value = 253631
$ws.Cells.Item(37, 1).Value = 'US Fund for UNICEF'
$ws.Cells.Item(37, 2).Value = 2454497
$ws.Cells.Item(38, 1).Value = 'UNICEF National Committee/France'
$ws.Cells.Item(38, 2).Value = 336700
$ws.Cells.Item(39, 1).Value = 'UNICEF National Committee/Portugal'
$ws.Cells.Item(39, 2).Value = 1129995
$ws.Cells.Item(40, 1).Value = 'UNICEF National Committee/Germany'
$ws.Cells.Item(40, 2).Value = 7175368
$ws.Cells.Item(41, 1).Value = 'Luxembourg, Government of'
$ws.Cells.Item(41, 2).Value = 278087
$ws.Cells.Item(42, 1).Value = 'Mozambique, Government of'
$ws.Cells.Item(42, 2).Value = 4493247
$ws.Cells.Item(43, 1).Value = 'UNICEF National Committee/Japan'
$ws.Cells.Item(43, 2).Value = 167052
$ws.Cells.Item(44, 1).Value = 'UNICEF National Committee/United Kingdom'
$ws.Cells.Item(44, 2).Value = 171977
$ws.Cells.Item(45, 1).Value = 'Estonia, Government of'
$ws.Cells.Item(45, 2).Value = 33076
$ws.Cells.Item(46, 1).Value = 'Italy, Government of'
$ws.Cells.Item(46, 2).Value = 5119678
$ws.Cells.Item(47, 1).Value = 'International Labour Organization'
$ws.Cells.Item(47, 2).Value = 355000
$ws.Cells.Item(48, 1).Value = 'World Food Programme'
$ws.Cells.Item(48, 2).Value = 2940725
$ws.Cells.Item(49, 1).Value = 'Monaco, Government of'
$ws.Cells.Item(49, 2).Value = 111483
$ws.Cells.Item(50, 1).Value = 'Private (individuals & organizations)'
$ws.Cells.Item(50, 2).Value = 1927645
$ws.Cells.Item(51, 1).Value = 'Malta, Government of'
$ws.Cells.Item(51, 2).Value = 22548
$ws.Cells.Item(52, 1).Value = 'Belgium, Government of'
$ws.Cells.Item(52, 2).Value = 3012377
$ws.Cells.Item(53, 1).Value = 'United Arab Emirates, Government of'
$ws.Cells.Item(53, 2).Value = 2180000
$ws.Cells.Item(54, 1).Value = 'International Organization for Migration'
$ws.Cells.Item(54, 2).Value = 497221
$ws.Cells.Item(55, 1).Value = 'China, Government of'
$ws.Cells.Item(55, 2).Value = 5062301
$ws.Cells.Item(56, 1).Value = 'World Bank'
$ws.Cells.Item(56, 2).Value = 10949258
$ws.Cells.Item(57, 1).Value = 'UN Delivering AS ONE'
$ws.Cells.Item(57, 2).Value = 647339
$ws.Cells.Item(58, 1).Value = 'Russian Federation, Government of'
$ws.Cells.Item(58, 2).Value = 1500000
$ws.Cells.Item(59, 1).Value = 'European Commission'
$ws.Cells.Item(59, 2).Value = 2688889
$ws.Cells.Item(60, 1).Value = 'Liechtenstein, Government of'
$ws.Cells.Item(60, 2).Value = 130775
$ws.Cells.Item(61, 1).Value = 'Jersey Overseas Aid'
$ws.Cells.Item(61, 2).Value = 190856
$ws.Cells.Item(62, 1).Value = 'India, Government of'
$ws.Cells.Item(62, 2).Value = 0
$ws.Cells.Item(63, 1).Value = 'Morocco, Government of'
$ws.Cells.Item(63, 2).Value = 0
$ws.Cells.Item(64, 1).Value = 'Turkey, Government of'
$ws.Cells.Item(64, 2).Value = 0
$ws.Cells.Item(65, 1).Value = 'Colombia, Government of'
$ws.Cells.Item(65, 2).Value = 50000
$ws.Cells.Item(66, 1).Value = 'Egypt, Government of'
$ws.Cells.Item(66, 2).Value = 237333
$ws.Cells.Item(67, 1).Value = 'Saudi Arabia (Kingdom of), Government of'
$ws.Cells.Item(67, 2).Value = 105333
$ws.Cells.Item(69, 1).Value = 'United Nations Population Fund'
$ws.Cells.Item(69, 2).Value = 475375
$ws.Cells.Item(70, 1).Value = 'New Zealand, Government of'
$ws.Cells.Item(70, 2).Value = 334672
$ws.Cells.Item(71, 1).Value = 'Azerbaijan, Government of'
$ws.Cells.Item(71, 2).Value = 100000
$ws.Cells.Item(72, 1).Value = 'Food & Agriculture Organization of the United Nations'
$ws.Cells.Item(72, 2).Value = 500000
$ws.Cells.Item(73, 1).Value = 'Netherlands, Government of'
$ws.Cells.Item(73, 2).Value = 3829322
$ws.Cells.Item(74, 1).Value = 'ACT Alliance / Church of Sweden'
$ws.Cells.Item(74, 2).Value = 898848
